$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update selection on "Binary Search 2" (it will no longer be the active
#    tab once "Notes" is added/activated below).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Binary Search 2")
$ws2.Activate()
$ws2.Range("B2:F2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) Add the new "Notes" worksheet after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Notes"

# Reuse the existing header formatting (style used by B2:F2 on the other
# sheets) so the shared style table indices line up instead of growing.
$ws2.Range("B2:F2").Copy()
$newSheet.Range("B2:F2").PasteSpecial(-4122)

# Header row (row 2)
$newSheet.Range("B2").Value = "S.no."
$newSheet.Range("C2").Value = "Page No. in notes"
$newSheet.Range("D2").Value = "Question"
$newSheet.Range("E2").Value = "Link"
$newSheet.Range("F2").Value = "Github Link"

# Data row (row 3)
$newSheet.Range("B3").Value = 1
$newSheet.Range("C3").Value = "Notes 1"
$newSheet.Range("D3").Value = "Find HCF or GCD"

# Hyperlink for the Link cell, then restore the same "Hyperlink" cell format
# used elsewhere in the workbook (F3 on "Binary Search 2").
$newSheet.Hyperlinks.Add($newSheet.Range("E3"), "https://www.scaler.com/academy/mentee-dashboard/class/28735/session") | Out-Null
$ws2.Range("F3").Copy()
$newSheet.Range("E3").PasteSpecial(-4122)

# Column widths matching the other sheets' look & feel.
$newSheet.Columns.Item(3).ColumnWidth = 14.5
$newSheet.Columns.Item(4).ColumnWidth = 13.5
$newSheet.Columns.Item(5).ColumnWidth = 23.666666666666668

# Row height for the wrapped data row.
$newSheet.Rows.Item(3).RowHeight = 57.6

# Make "Notes" the active sheet/tab with F3 selected.
$newSheet.Activate()
$newSheet.Range("F3").Select() | Out-Null

$excel.CutCopyMode = 0
